$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.053.77'
$ws.Range('E2').Value = '  +4.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.651.20'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.94'
$ws.Range('E5').Value = '  +7.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.22'
$ws.Range('E6').Value = '  +3.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.608'
$ws.Range('E8').Value = '  +7.21%  '
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('E10').Value = '  +4.68%  '
$ws.Range('E11').Value = '  +6.61%  '
$ws.Range('E12').Value = '  +3.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.119.68'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '60.926.63'
$ws.Range('E14').Value = '  +4.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.87'
$ws.Range('E15').Value = '  +6.15%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000138'
$ws.Range('E16').Value = '  +5.26%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.671.95'
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '344.98'
$ws.Range('E19').Value = '  +3.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.49'
$ws.Range('E20').Value = '  +3.90%  '
$ws.Range('E21').Value = '  +3.53%  '
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.97'
$ws.Range('E24').Value = '  +1.06%  '
$ws.Range('E25').Value = '  +7.39%  '
$ws.Range('E26').Value = '  +1.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.991'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('E28').Value = '  +5.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0790'
$ws.Range('E29').Value = '  +7.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  +5.34%  '
$ws.Range('E32').Value = '  +5.03%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.29'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '154.88'
$ws.Range('E34').Value = '  +3.10%  '
$ws.Range('E35').Value = '  +6.27%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.19'
$ws.Range('E36').Value = '  +8.68%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.913'
$ws.Range('E37').Value = '  +6.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.917'
$ws.Range('E38').Value = '  +12.95%  '
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('E40').Value = '  +8.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '308.13'
$ws.Range('E41').Value = '  +10.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.67'
$ws.Range('E42').Value = '  +3.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.997'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0983'
$ws.Range('E45').Value = '  +5.26%  '
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.48'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.39'
$ws.Range('E48').Value = '  +12.92%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.93'
$ws.Range('E51').Value = '  +5.98%  '
